$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3557.3809
$ws.Range("I69").Value = 2844.1667
$ws.Range("J69").Value = 4508.3335
$ws.Range("K69").Value = 8532.500100000001
$ws.Range("L69").Value = 13525.0005
$ws.Range("M69").Value = -7658.500100000001
$ws.Range("N69").Value = -15273.0005
$ws.Range("H72").Value = 3557.3809
$ws.Range("I72").Value = 2844.1667
$ws.Range("J72").Value = 4508.3335
$ws.Range("K72").Value = 25597.5003
$ws.Range("L72").Value = 40575.0015
$ws.Range("M72").Value = -21229.5003
$ws.Range("N72").Value = -49311.0015
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H129").Value = 1018.05884
$ws.Range("J129").Value = 1128.3489
$ws.Range("L129").Value = 3385.0467
$ws.Range("N129").Value = -13385.0467
$ws.Range("H132").Value = 1287.3846
$ws.Range("I132").Value = 1262.0938
$ws.Range("K132").Value = 3786.2814
$ws.Range("M132").Value = -1256.2814
$ws.Range("H137").Value = 1711.5428
$ws.Range("I137").Value = 1103.5
$ws.Range("K137").Value = 3310.5
$ws.Range("M137").Value = -760.5
$ws.Range("H138").Value = 1708.2222
$ws.Range("I138").Value = 1325.0889
$ws.Range("J138").Value = 3623.889
$ws.Range("K138").Value = 3975.2667
$ws.Range("L138").Value = 10871.667
$ws.Range("M138").Value = 1164.7333
$ws.Range("N138").Value = -21151.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4251.8105
$ws.Range("I32").Value = 3002.4597
$ws.Range("J32").Value = 17838.5
$ws.Range("K32").Value = 3002.4597
$ws.Range("L32").Value = 17838.5
$ws.Range("M32").Value = -2715.4597
$ws.Range("N32").Value = -18412.5
$ws.Range("H74").Value = 4253
$ws.Range("I74").Value = 4253
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4253
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3379
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 4253
$ws.Range("I77").Value = 4253
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 21265
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -16897
$ws.Range("N77").ClearContents()
$ws.Range("H132").Value = 2758.0344
$ws.Range("I132").Value = 1516.7059
$ws.Range("J132").Value = 4516.5835
$ws.Range("K132").Value = 4550.1177
$ws.Range("L132").Value = 13549.7505
$ws.Range("M132").Value = -2020.1177
$ws.Range("N132").Value = -18609.7505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1944.7142
$ws.Range("I20").Value = 1979.6666
$ws.Range("J20").Value = 1881.8
$ws.Range("K20").Value = 1979.6666
$ws.Range("L20").Value = 1881.8
$ws.Range("M20").Value = -1732.6666
$ws.Range("N20").Value = -2375.8
$ws.Range("H105").Value = 5684130
$ws.Range("I105").Value = 5684130
$ws.Range("K105").Value = 5684130
$ws.Range("M105").Value = -5682383
$ws.Range("H134").Value = 1686.7587
$ws.Range("I134").Value = 1401.5
$ws.Range("J134").Value = 3056
$ws.Range("K134").Value = 4204.5
$ws.Range("L134").Value = 9168
$ws.Range("M134").Value = -1669.5
$ws.Range("N134").Value = -14238

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1030.9524
$ws.Range("I105").Value = 1030.9524
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1030.9524
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 716.0476000000001
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 2211.842
$ws.Range("I132").Value = 943.9091
$ws.Range("J132").Value = 3955.25
$ws.Range("K132").Value = 2831.7273
$ws.Range("L132").Value = 11865.75
$ws.Range("M132").Value = -301.7273
$ws.Range("N132").Value = -16925.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5650.909
$ws.Range("I131").Value = 456.36365
$ws.Range("J131").Value = 10845.454
$ws.Range("K131").Value = 1369.09095
$ws.Range("L131").Value = 32536.362
$ws.Range("M131").Value = 3670.90905
$ws.Range("N131").Value = -42616.362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5320.0933
$ws.Range("I70").Value = 4666.1
$ws.Range("J70").Value = 6829.3076
$ws.Range("K70").Value = 4666.1
$ws.Range("L70").Value = 6829.3076
$ws.Range("M70").Value = -4396.1
$ws.Range("N70").Value = -7369.3076
$ws.Range("H73").Value = 5320.0933
$ws.Range("I73").Value = 4666.1
$ws.Range("J73").Value = 6829.3076
$ws.Range("K73").Value = 4666.1
$ws.Range("L73").Value = 6829.3076
$ws.Range("M73").Value = -3730.1
$ws.Range("N73").Value = -8701.3076
$ws.Range("H80").Value = 3998.3333
$ws.Range("I80").Value = 3997.8
$ws.Range("J80").Value = 3999.4
$ws.Range("K80").Value = 3997.8
$ws.Range("L80").Value = 3999.4
$ws.Range("M80").Value = -2999.8
$ws.Range("N80").Value = -5995.4
$ws.Range("H83").Value = 3998.3333
$ws.Range("I83").Value = 3997.8
$ws.Range("J83").Value = 3999.4
$ws.Range("K83").Value = 19989
$ws.Range("L83").Value = 19997
$ws.Range("M83").Value = -14997
$ws.Range("N83").Value = -29981
$ws.Range("H126").Value = 1756.3158
$ws.Range("I126").Value = 1521.091
$ws.Range("J126").Value = 3308.8
$ws.Range("K126").Value = 4563.272999999999
$ws.Range("L126").Value = 9926.400000000001
$ws.Range("M126").Value = -2093.272999999999
$ws.Range("N126").Value = -14866.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7222.4526
$ws.Range("I132").Value = 6468.619
$ws.Range("J132").Value = 10100.728
$ws.Range("K132").Value = 19405.857
$ws.Range("L132").Value = 30302.184
$ws.Range("M132").Value = -16875.857
$ws.Range("N132").Value = -35362.18399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3488.889
$ws.Range("J62").Value = 3488.889
$ws.Range("L62").Value = 3488.889
$ws.Range("N62").Value = -4736.889
$ws.Range("H65").Value = 3488.889
$ws.Range("J65").Value = 3488.889
$ws.Range("L65").Value = 17444.445
$ws.Range("N65").Value = -23684.445
$ws.Range("H81").Value = 1324.4117
$ws.Range("I81").Value = 800.25
$ws.Range("K81").Value = 1600.5
$ws.Range("M81").Value = -539.5
$ws.Range("H84").Value = 1324.4117
$ws.Range("I84").Value = 800.25
$ws.Range("K84").Value = 8002.5
$ws.Range("M84").Value = -2698.5
$ws.Range("H132").Value = 1265.9818
$ws.Range("I132").Value = 924.5
$ws.Range("J132").Value = 3011.3333
$ws.Range("K132").Value = 2773.5
$ws.Range("L132").Value = 9033.999899999999
$ws.Range("M132").Value = -243.5
$ws.Range("N132").Value = -14093.9999
$ws.Range("H136").Value = 6291276
$ws.Range("I136").Value = 7937507.5
$ws.Range("J136").Value = 5664.091
$ws.Range("K136").Value = 23812522.5
$ws.Range("L136").Value = 16992.273
$ws.Range("M136").Value = -23809972.5
$ws.Range("N136").Value = -22092.273
